$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Copy($ws.Range("F1"))
$ws.Range("F1").Value = "time_taken"

$timestamps = @(
    "2021-10-05 10:52:45.844649",
    "2021-10-05 10:52:45.844660",
    "2021-10-05 10:52:45.844664",
    "2021-10-05 10:52:45.844667",
    "2021-10-05 10:52:45.844671",
    "2021-10-05 10:52:45.844674",
    "2021-10-05 10:52:45.844677",
    "2021-10-05 10:52:45.844680",
    "2021-10-05 10:52:45.844683",
    "2021-10-05 10:52:45.844686",
    "2021-10-05 10:52:45.844689",
    "2021-10-05 10:52:45.844692",
    "2021-10-05 10:52:45.844695"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timestamps[$i]
}
